{"js": "// Update the two-digit multiplication answer table (20 rows x 5 cols)\n// with the new set of equations, cell by cell, in row-major (reading) order.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst newValues = [\n  [\"39\u00d749=1911\", \"48\u00d727=1296\", \"10\u00d748=480\", \"58\u00d726=1508\", \"30\u00d759=1770\"],\n  [\"29\u00d7100=2900\", \"40\u00d760=2400\", \"69\u00d743=2967\", \"47\u00d746=2162\", \"49\u00d750=2450\"],\n  [\"34\u00d715=510\", \"43\u00d769=2967\", \"12\u00d761=732\", \"78\u00d736=2808\", \"59\u00d752=3068\"],\n  [\"41\u00d765=2665\", \"76\u00d746=3496\", \"51\u00d711=561\", \"13\u00d797=1261\", \"56\u00d757=3192\"],\n  [\"88\u00d748=4224\", \"29\u00d740=1160\", \"30\u00d711=330\", \"18\u00d725=450\", \"61\u00d765=3965\"],\n  [\"81\u00d765=5265\", \"75\u00d790=6750\", \"49\u00d730=1470\", \"54\u00d745=2430\", \"17\u00d792=1564\"],\n  [\"23\u00d757=1311\", \"47\u00d782=3854\", \"100\u00d743=4300\", \"22\u00d771=1562\", \"70\u00d777=5390\"],\n  [\"34\u00d742=1428\", \"36\u00d760=2160\", \"93\u00d768=6324\", \"95\u00d763=5985\", \"16\u00d795=1520\"],\n  [\"76\u00d765=4940\", \"68\u00d753=3604\", \"61\u00d750=3050\", \"39\u00d771=2769\", \"11\u00d757=627\"],\n  [\"10\u00d750=500\", \"69\u00d798=6762\", \"48\u00d772=3456\", \"72\u00d717=1224\", \"85\u00d711=935\"],\n  [\"17\u00d748=816\", \"38\u00d739=1482\", \"82\u00d759=4838\", \"29\u00d776=2204\", \"30\u00d781=2430\"],\n  [\"57\u00d750=2850\", \"11\u00d783=913\", \"48\u00d787=4176\", \"82\u00d734=2788\", \"95\u00d790=8550\"],\n  [\"77\u00d728=2156\", \"50\u00d728=1400\", \"98\u00d761=5978\", \"12\u00d798=1176\", \"64\u00d758=3712\"],\n  [\"84\u00d790=7560\", \"20\u00d778=1560\", \"13\u00d789=1157\", \"47\u00d728=1316\", \"86\u00d792=7912\"],\n  [\"63\u00d790=5670\", \"44\u00d736=1584\", \"22\u00d719=418\", \"41\u00d751=2091\", \"41\u00d771=2911\"],\n  [\"19\u00d795=1805\", \"54\u00d711=594\", \"87\u00d776=6612\", \"90\u00d768=6120\", \"80\u00d749=3920\"],\n  [\"35\u00d759=2065\", \"19\u00d734=646\", \"66\u00d783=5478\", \"55\u00d772=3960\", \"11\u00d768=748\"],\n  [\"99\u00d778=7722\", \"46\u00d755=2530\", \"57\u00d780=4560\", \"98\u00d755=5390\", \"85\u00d724=2040\"],\n  [\"30\u00d732=960\", \"47\u00d768=3196\", \"63\u00d721=1323\", \"91\u00d734=3094\", \"65\u00d791=5915\"],\n  [\"95\u00d780=7600\", \"12\u00d775=900\", \"89\u00d797=8633\", \"52\u00d714=728\", \"36\u00d731=1116\"],\n];\n\nif (table.rowCount !== newValues.length) {\n  throw new Error(\n    `Table row count (${table.rowCount}) does not match replacement data (${newValues.length} rows)`\n  );\n}\n\n// Assigning the full 2-D array replaces each cell's text in place while\n// preserving existing paragraph/run formatting (font, size, alignment, etc).\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Update the two-digit multiplication answer table (20 rows x 5 cols)\n# with the new set of equations, cell by cell, in row-major order.\n$newValues = @(\n    \"39\u00d749=1911\", \"48\u00d727=1296\", \"10\u00d748=480\", \"58\u00d726=1508\", \"30\u00d759=1770\",\n    \"29\u00d7100=2900\", \"40\u00d760=2400\", \"69\u00d743=2967\", \"47\u00d746=2162\", \"49\u00d750=2450\",\n    \"34\u00d715=510\", \"43\u00d769=2967\", \"12\u00d761=732\", \"78\u00d736=2808\", \"59\u00d752=3068\",\n    \"41\u00d765=2665\", \"76\u00d746=3496\", \"51\u00d711=561\", \"13\u00d797=1261\", \"56\u00d757=3192\",\n    \"88\u00d748=4224\", \"29\u00d740=1160\", \"30\u00d711=330\", \"18\u00d725=450\", \"61\u00d765=3965\",\n    \"81\u00d765=5265\", \"75\u00d790=6750\", \"49\u00d730=1470\", \"54\u00d745=2430\", \"17\u00d792=1564\",\n    \"23\u00d757=1311\", \"47\u00d782=3854\", \"100\u00d743=4300\", \"22\u00d771=1562\", \"70\u00d777=5390\",\n    \"34\u00d742=1428\", \"36\u00d760=2160\", \"93\u00d768=6324\", \"95\u00d763=5985\", \"16\u00d795=1520\",\n    \"76\u00d765=4940\", \"68\u00d753=3604\", \"61\u00d750=3050\", \"39\u00d771=2769\", \"11\u00d757=627\",\n    \"10\u00d750=500\", \"69\u00d798=6762\", \"48\u00d772=3456\", \"72\u00d717=1224\", \"85\u00d711=935\",\n    \"17\u00d748=816\", \"38\u00d739=1482\", \"82\u00d759=4838\", \"29\u00d776=2204\", \"30\u00d781=2430\",\n    \"57\u00d750=2850\", \"11\u00d783=913\", \"48\u00d787=4176\", \"82\u00d734=2788\", \"95\u00d790=8550\",\n    \"77\u00d728=2156\", \"50\u00d728=1400\", \"98\u00d761=5978\", \"12\u00d798=1176\", \"64\u00d758=3712\",\n    \"84\u00d790=7560\", \"20\u00d778=1560\", \"13\u00d789=1157\", \"47\u00d728=1316\", \"86\u00d792=7912\",\n    \"63\u00d790=5670\", \"44\u00d736=1584\", \"22\u00d719=418\", \"41\u00d751=2091\", \"41\u00d771=2911\",\n    \"19\u00d795=1805\", \"54\u00d711=594\", \"87\u00d776=6612\", \"90\u00d768=6120\", \"80\u00d749=3920\",\n    \"35\u00d759=2065\", \"19\u00d734=646\", \"66\u00d783=5478\", \"55\u00d772=3960\", \"11\u00d768=748\",\n    \"99\u00d778=7722\", \"46\u00d755=2530\", \"57\u00d780=4560\", \"98\u00d755=5390\", \"85\u00d724=2040\",\n    \"30\u00d732=960\", \"47\u00d768=3196\", \"63\u00d721=1323\", \"91\u00d734=3094\", \"65\u00d791=5915\",\n    \"95\u00d780=7600\", \"12\u00d775=900\", \"89\u00d797=8633\", \"52\u00d714=728\", \"36\u00d731=1116\"\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$rowCount = $tbl.Rows.Count\n$colCount = $tbl.Columns.Count\nif ($rowCount * $colCount -ne $newValues.Count) {\n    throw \"Table shape ($rowCount x $colCount) does not match replacement data ($($newValues.Count) values)\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $cell.Range.Text = $newValues[$i]\n        $i++\n    }\n}\n"}
